$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.003" or
# "332.34" are not silently reinterpreted as numbers (which would drop
# trailing zeros / change formatting). We restore the default "Normal"
# style afterwards so no stray style index is left on the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Cell value updates, row by row, matching the source diff ---

# Row 2
$ws.Range("D2").Value = '30.484.37'
$ws.Range("E2").Value = '  -0.94%  '

# Row 3
$ws.Range("D3").Value = '2.103.12'
$ws.Range("E3").Value = '  -0.39%  '

# Row 4
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.14%  '

# Row 5
$ws.Range("D5").Value = '332.34'
$ws.Range("E5").Value = '  -0.33%  '

# Row 6
$ws.Range("E6").Value = '  +0.17%  '

# Row 7
$ws.Range("D7").Value = '0.5228'
$ws.Range("E7").Value = '  -1.21%  '

# Row 8
$ws.Range("D8").Value = '0.4494'
$ws.Range("E8").Value = '  +2.02%  '

# Row 9
$ws.Range("D9").Value = '53.78'
$ws.Range("E9").Value = '  +16.97%  '

# Row 10
$ws.Range("D10").Value = '0.08922'
$ws.Range("E10").Value = '  -1.03%  '

# Row 11
$ws.Range("D11").Value = '1.157'
$ws.Range("E11").Value = '  -1.89%  '

# Row 12
$ws.Range("E12").Value = '  -2.64%  '

# Row 13
$ws.Range("D13").Value = '2.094.87'
$ws.Range("E13").Value = '  -0.57%  '

# Row 14
$ws.Range("D14").Value = '6.743'
$ws.Range("E14").Value = '  -0.13%  '

# Row 15
$ws.Range("D15").Value = '7.724'
$ws.Range("E15").Value = '  -0.94%  '

# Row 16
$ws.Range("D16").Value = '96.40'
$ws.Range("E16").Value = '  -1.10%  '

# Row 17
$ws.Range("E17").Value = '  +0.23%  '

# Row 18
$ws.Range("D18").Value = '0.00001125'
$ws.Range("E18").Value = '  -0.15%  '

# Row 19
$ws.Range("D19").Value = '0.06622'
$ws.Range("E19").Value = '  -0.65%  '

# Row 20
$ws.Range("E20").Value = '  +0.40%  '

# Row 21
$ws.Range("E21").Value = '  +0.08%  '

# Row 22
$ws.Range("D22").Value = '6.289'
$ws.Range("E22").Value = '  -1.12%  '

# Row 23
$ws.Range("D23").Value = '30.539.14'
$ws.Range("E23").Value = '  -0.93%  '

# Row 24
$ws.Range("D24").Value = '12.33'
$ws.Range("E24").Value = '  +0.65%  '

# Row 25
$ws.Range("D25").Value = '2.319'
$ws.Range("E25").Value = '  +2.72%  '

# Row 26
$ws.Range("D26").Value = '2.340.22'
$ws.Range("E26").Value = '  -0.70%  '

# Row 27
$ws.Range("D27").Value = '22.30'
$ws.Range("E27").Value = '  -2.35%  '

# Row 28
$ws.Range("D28").Value = '2.583'
$ws.Range("E28").Value = '  +0.97%  '

# Row 29
$ws.Range("D29").Value = '163.99'
$ws.Range("E29").Value = '  +0.93%  '

# Row 30
$ws.Range("D30").Value = '132.43'
$ws.Range("E30").Value = '  -0.49%  '

# Row 31
$ws.Range("D31").Value = '1.200'
$ws.Range("E31").Value = '  +2.52%  '

# Row 32
$ws.Range("E32").Value = '  -0.09%  '

# Row 33
$ws.Range("D33").Value = '1.666'
$ws.Range("E33").Value = '  +7.78%  '

# Row 34
$ws.Range("D34").Value = '6.151'
$ws.Range("E34").Value = '  -1.11%  '

# Row 35
$ws.Range("D35").Value = '3.943'
$ws.Range("E35").Value = '  -2.05%  '

# Row 36
$ws.Range("D36").Value = '10.45'
$ws.Range("E36").Value = '  +9.87%  '

# Row 37
$ws.Range("D37").Value = '0.02569'
$ws.Range("E37").Value = '  -1.31%  '

# Row 38
$ws.Range("D38").Value = '0.06779'
$ws.Range("E38").Value = '  +0.58%  '

# Row 39
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = '12.80'
$ws.Range("E39").Value = '  +1.00%  '

# Row 40
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").Value = '5.484'
$ws.Range("E40").Value = '  -0.86%  '

# Row 41
$ws.Range("D41").Value = '0.2267'
$ws.Range("E41").Value = '  -0.04%  '

# Row 42
$ws.Range("D42").Value = '0.6929'
$ws.Range("E42").Value = '  +1.16%  '

# Row 43
$ws.Range("D43").Value = '1.258'
$ws.Range("E43").Value = '  +1.34%  '

# Row 44
$ws.Range("D44").Value = '1.003'
$ws.Range("E44").Value = '  +0.16%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '14.00'
$ws.Range("E45").Value = '  -0.92%  '

# Row 46
$ws.Range("D46").Value = '0.6372'
$ws.Range("E46").Value = '  -1.28%  '

# Row 47
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '2.299'
$ws.Range("E47").Value = '  +3.04%  '

# Row 48
$ws.Range("D48").Value = '3.641'
$ws.Range("E48").Value = '  -0.73%  '

# Row 49
$ws.Range("D49").Value = '1.245'
$ws.Range("E49").Value = '  -2.37%  '

# Row 50
$ws.Range("D50").Value = '1.223'
$ws.Range("E50").Value = '  +5.58%  '

# Row 51
$ws.Range("D51").Value = '83.01'
$ws.Range("E51").Value = '  +0.79%  '

# Restore default styling on the price column so cells keep their text
# type without retaining an extra/unused number-format style index.
$priceRange.Style = "Normal"
